$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 18 through 28 (second batch of P08_4 samples), shrinking the
# used range from A1:P28 down to A1:P17.
$ws.Range("A18:P28").EntireRow.Delete()
